$d = $word.ActiveDocument

$found = $d.Content.Find.Execute("550×9=", $false, $false, $false, $false, $false, $true, 1, $false, "607×4=", 2)
if (-not $found) { throw "Could not find text: 550×9=" }
$found = $d.Content.Find.Execute("612×7=", $false, $false, $false, $false, $false, $true, 1, $false, "680×8=", 2)
if (-not $found) { throw "Could not find text: 612×7=" }
$found = $d.Content.Find.Execute("194×5=", $false, $false, $false, $false, $false, $true, 1, $false, "135×9=", 2)
if (-not $found) { throw "Could not find text: 194×5=" }
$found = $d.Content.Find.Execute("232×6=", $false, $false, $false, $false, $false, $true, 1, $false, "301×7=", 2)
if (-not $found) { throw "Could not find text: 232×6=" }
$found = $d.Content.Find.Execute("465×5=", $false, $false, $false, $false, $false, $true, 1, $false, "628×9=", 2)
if (-not $found) { throw "Could not find text: 465×5=" }
$found = $d.Content.Find.Execute("533×7=", $false, $false, $false, $false, $false, $true, 1, $false, "129×8=", 2)
if (-not $found) { throw "Could not find text: 533×7=" }
$found = $d.Content.Find.Execute("662×5=", $false, $false, $false, $false, $false, $true, 1, $false, "303×3=", 2)
if (-not $found) { throw "Could not find text: 662×5=" }
$found = $d.Content.Find.Execute("737×7=", $false, $false, $false, $false, $false, $true, 1, $false, "326×5=", 2)
if (-not $found) { throw "Could not find text: 737×7=" }
$found = $d.Content.Find.Execute("107×5=", $false, $false, $false, $false, $false, $true, 1, $false, "792×2=", 2)
if (-not $found) { throw "Could not find text: 107×5=" }
$found = $d.Content.Find.Execute("359×5=", $false, $false, $false, $false, $false, $true, 1, $false, "567×6=", 2)
if (-not $found) { throw "Could not find text: 359×5=" }
$found = $d.Content.Find.Execute("475×7=", $false, $false, $false, $false, $false, $true, 1, $false, "349×8=", 2)
if (-not $found) { throw "Could not find text: 475×7=" }
$found = $d.Content.Find.Execute("554×2=", $false, $false, $false, $false, $false, $true, 1, $false, "867×8=", 2)
if (-not $found) { throw "Could not find text: 554×2=" }
$found = $d.Content.Find.Execute("265×7=", $false, $false, $false, $false, $false, $true, 1, $false, "765×4=", 2)
if (-not $found) { throw "Could not find text: 265×7=" }
$found = $d.Content.Find.Execute("738×4=", $false, $false, $false, $false, $false, $true, 1, $false, "823×8=", 2)
if (-not $found) { throw "Could not find text: 738×4=" }
$found = $d.Content.Find.Execute("636×9=", $false, $false, $false, $false, $false, $true, 1, $false, "890×5=", 2)
if (-not $found) { throw "Could not find text: 636×9=" }
$found = $d.Content.Find.Execute("393×6=", $false, $false, $false, $false, $false, $true, 1, $false, "107×5=", 2)
if (-not $found) { throw "Could not find text: 393×6=" }
$found = $d.Content.Find.Execute("700×3=", $false, $false, $false, $false, $false, $true, 1, $false, "491×2=", 2)
if (-not $found) { throw "Could not find text: 700×3=" }
$found = $d.Content.Find.Execute("241×6=", $false, $false, $false, $false, $false, $true, 1, $false, "814×8=", 2)
if (-not $found) { throw "Could not find text: 241×6=" }
$found = $d.Content.Find.Execute("124×8=", $false, $false, $false, $false, $false, $true, 1, $false, "513×3=", 2)
if (-not $found) { throw "Could not find text: 124×8=" }
$found = $d.Content.Find.Execute("321×4=", $false, $false, $false, $false, $false, $true, 1, $false, "760×6=", 2)
if (-not $found) { throw "Could not find text: 321×4=" }
$found = $d.Content.Find.Execute("418×2=", $false, $false, $false, $false, $false, $true, 1, $false, "367×7=", 2)
if (-not $found) { throw "Could not find text: 418×2=" }
$found = $d.Content.Find.Execute("934×5=", $false, $false, $false, $false, $false, $true, 1, $false, "806×4=", 2)
if (-not $found) { throw "Could not find text: 934×5=" }
$found = $d.Content.Find.Execute("416×3=", $false, $false, $false, $false, $false, $true, 1, $false, "565×3=", 2)
if (-not $found) { throw "Could not find text: 416×3=" }
$found = $d.Content.Find.Execute("365×6=", $false, $false, $false, $false, $false, $true, 1, $false, "802×3=", 2)
if (-not $found) { throw "Could not find text: 365×6=" }
$found = $d.Content.Find.Execute("829×7=", $false, $false, $false, $false, $false, $true, 1, $false, "442×5=", 2)
if (-not $found) { throw "Could not find text: 829×7=" }
